# "New Code with API" — rename sheets, refresh test-data rows, add new
# guest-checkout / zip-code rows, and restore the UI selection / active-tab
# state recorded by Excel when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheets (Products -> PDP, checkOut -> Checkout)
# ---------------------------------------------------------------------
$wsPDP = $wb.Worksheets.Item("Products")
$wsPDP.Name = "PDP"

$wsCheckout = $wb.Worksheets.Item("checkOut")
$wsCheckout.Name = "Checkout"

# ---------------------------------------------------------------------
# 2. PDP sheet: product under test swapped + key renames
# ---------------------------------------------------------------------
$wsPDP.Range("B2").Value = "TimeWise Miracle Set 3D®"
$wsPDP.Range("A3").Value = "Product"
$wsPDP.Range("A6").Value = "productCount"

# ---------------------------------------------------------------------
# 3. Checkout sheet: key renames + two new rows (guest checkout / button)
# ---------------------------------------------------------------------
$wsCheckout.Range("A2").Value = "ZipCode"
$wsCheckout.Range("A3").Value = "guestCheckOutHeading"
$wsCheckout.Range("A4").Value = "guestChek-Out"
$wsCheckout.Range("B4").Value = "Guest Checkout"
$wsCheckout.Range("A5").Value = "checkOutbtn"
$wsCheckout.Range("B5").Value = "Checkout"

# ---------------------------------------------------------------------
# 4. ShipperForm sheet: new ZipCode row, copied (value + format) from the
#    now-updated Checkout!A2:B2 row
# ---------------------------------------------------------------------
$wsShipper = $wb.Worksheets.Item("ShipperForm")
$wsCheckout.Range("A2:B2").Copy($wsShipper.Range("A10:B10"))

# ---------------------------------------------------------------------
# 5. Restore each sheet's remembered selection
# ---------------------------------------------------------------------
$wsSetUp = $wb.Worksheets.Item("SetUp")
$wsSetUp.Range("B2").Select()

$wsSearchPage = $wb.Worksheets.Item("SearchPage")
$wsSearchPage.Range("A2").Select()

$wsPDP.Range("A6").Select()

$wsCheckout.Range("B5").Select()

$wsShipper.Range("A10:B10").Select()

# ---------------------------------------------------------------------
# 6. Checkout becomes the active tab (was SearchPage before)
# ---------------------------------------------------------------------
$wsCheckout.Activate()
